$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "PT Borwita" lease record (row 5). Excel shifts all
# subsequent rows up by one and the used range shrinks accordingly.
$ws.Rows.Item(5).Delete()
